# Generate Report for Archive
#
# 1) Status text changed from "Ready for handoff" to "In Translation"
#    (this shared string is used on the Overview sheet in E2/F2, and on
#    each language sheet's "Status" column, cell C2).
# 2) The two status columns (zh-cn/de-de) on the Overview sheet, and the
#    Status column on each language sheet, are narrowed to match the new,
#    shorter status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status value everywhere it appears ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Narrow the status columns to fit the new text ---
# (Target author width is ~13.41 characters; this runtime's ColumnWidth
# setter, like Excel's, only keeps whole-pixel precision, so 12.5 is the
# input that lands on the closest representable stored width.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
